# Insert a new weekly price-report row for "Poroto verde" above the
# existing row 26 (Vega Modelo de Temuco / La Araucania), pushing the
# previously-existing rows 26-84 down to 27-85.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 26; Excel shifts rows 26:84 -> 27:85 and
# grows the used range / dimension automatically.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new record.
$ws.Cells.Item(26, 1).Value  = 10
$ws.Cells.Item(26, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(26, 3).Value  = "La Araucanía"
$ws.Cells.Item(26, 4).Value  = 44497
$ws.Cells.Item(26, 5).Value  = 9
$ws.Cells.Item(26, 6).Value  = 100112031
$ws.Cells.Item(26, 7).Value  = "Poroto verde"
$ws.Cells.Item(26, 8).Value  = "Sin especificar"
$ws.Cells.Item(26, 9).Value  = "Primera"
$ws.Cells.Item(26, 10).Value = 50
$ws.Cells.Item(26, 11).Value = 45000
$ws.Cells.Item(26, 12).Value = 45000
$ws.Cells.Item(26, 13).Value = 45000
$ws.Cells.Item(26, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(26, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(26, 16).Value = 1800
$ws.Cells.Item(26, 17).Value = 25
$ws.Cells.Item(26, 18).Value = "Hortaliza"

# Keep the date column formatted like the rest of column D.
$ws.Cells.Item(26, 4).NumberFormat = $ws.Cells.Item(27, 4).NumberFormat
